$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(401,9,48,67,75,45),
    @(301,6,45,30,60,45),
    @(501,9,52,30,75,45),
    @(701,3,90,45,97,15),
    @(1201,2,10,10,10,10),
    @(1202,2,10,10,10,10),
    @(902,1,0,0,0,0),
    @(201,9,30,15,45,30),
    @(801,3,67,65,52,45),
    @(1203,3,15,15,15,15),
    @(101,9,30,15,60,15),
    @(901,16,15,45,60,60),
    @(1001,18,30,75,60,72),
    @(601,9,60,67,60,42),
    @(502,0,4,0,0,0),
    @(1101,0,15,30,30,0),
    @(3,0,3,3,3,3),
    @(2,0,2,2,2,2),
    @(1,0,2,2,2,2),
    @(802,0,4,5,4,0),
    @(402,0,0,4,0,0),
    @(602,0,0,4,0,9),
    @(702,0,0,0,4,0),
    @(1002,0,0,0,0,9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col + 1).Value = $rowValues[$col]
    }
}
